$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.24254846572876
$ws.Range("B1").Value = 1.314867258071899
$ws.Range("C1").Value = 1.488025426864624
$ws.Range("D1").Value = 2.330355882644653
$ws.Range("E1").Value = 15
